# Apply edits described by the diff: update Fecha/Volumen/prices/Unidad/Origen
# values for rows 3-19, and append two new data rows (20, 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 4).Value = 44425
$ws.Cells.Item(3, 10).Value = 10
$ws.Cells.Item(3, 15).Value = 'Región de Arica y Parinacota'

# Row 4
$ws.Cells.Item(4, 4).Value = 44329
$ws.Cells.Item(4, 11).Value = 20000
$ws.Cells.Item(4, 12).Value = 20000
$ws.Cells.Item(4, 13).Value = 20000
$ws.Cells.Item(4, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(4, 16).Value = 1333
$ws.Cells.Item(4, 17).Value = 15

# Row 5
$ws.Cells.Item(5, 4).Value = 44424
$ws.Cells.Item(5, 10).Value = 30

# Row 6
$ws.Cells.Item(6, 4).Value = 44340
$ws.Cells.Item(6, 10).Value = 40
$ws.Cells.Item(6, 11).Value = 18000
$ws.Cells.Item(6, 12).Value = 18000
$ws.Cells.Item(6, 13).Value = 18000
$ws.Cells.Item(6, 15).Value = 'Perú'
$ws.Cells.Item(6, 16).Value = 900

# Row 7
$ws.Cells.Item(7, 4).Value = 44175
$ws.Cells.Item(7, 10).Value = 20
$ws.Cells.Item(7, 11).Value = 20000
$ws.Cells.Item(7, 12).Value = 20000
$ws.Cells.Item(7, 13).Value = 20000
$ws.Cells.Item(7, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(7, 16).Value = 1333

# Row 8
$ws.Cells.Item(8, 4).Value = 44385
$ws.Cells.Item(8, 10).Value = 18
$ws.Cells.Item(8, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(8, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(8, 16).Value = 1000
$ws.Cells.Item(8, 17).Value = 20

# Row 9
$ws.Cells.Item(9, 4).Value = 44321
$ws.Cells.Item(9, 10).Value = 15
$ws.Cells.Item(9, 11).Value = 25000
$ws.Cells.Item(9, 12).Value = 25000
$ws.Cells.Item(9, 13).Value = 25000
$ws.Cells.Item(9, 15).Value = 'Perú'
$ws.Cells.Item(9, 16).Value = 1667

# Row 10
$ws.Cells.Item(10, 4).Value = 44364
$ws.Cells.Item(10, 10).Value = 15
$ws.Cells.Item(10, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(10, 15).Value = 'Perú'
$ws.Cells.Item(10, 16).Value = 1333
$ws.Cells.Item(10, 17).Value = 15

# Row 11
$ws.Cells.Item(11, 4).Value = 44315
$ws.Cells.Item(11, 10).Value = 30

# Row 12
$ws.Cells.Item(12, 4).Value = 44315
$ws.Cells.Item(12, 10).Value = 30
$ws.Cells.Item(12, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(12, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(12, 16).Value = 1000
$ws.Cells.Item(12, 17).Value = 20

# Row 13
$ws.Cells.Item(13, 4).Value = 44389
$ws.Cells.Item(13, 10).Value = 45

# Row 14
$ws.Cells.Item(14, 4).Value = 44294
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(14, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(14, 15).Value = 'Perú'
$ws.Cells.Item(14, 16).Value = 1333
$ws.Cells.Item(14, 17).Value = 15

# Row 15
$ws.Cells.Item(15, 4).Value = 44369

# Row 16
$ws.Cells.Item(16, 4).Value = 44369
$ws.Cells.Item(16, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(16, 16).Value = 1000
$ws.Cells.Item(16, 17).Value = 20

# Row 17
$ws.Cells.Item(17, 4).Value = 44161

# Row 18
$ws.Cells.Item(18, 4).Value = 44188

# Row 19
$ws.Cells.Item(19, 4).Value = 44316

# New row 20
$ws.Cells.Item(20, 1).Value = 10
$ws.Cells.Item(20, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(20, 3).Value = 'La Araucanía'
$ws.Cells.Item(20, 4).Value = 44186
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(20, 6).Value = 100114002
$ws.Cells.Item(20, 7).Value = 'Camote'
$ws.Cells.Item(20, 8).Value = 'Sin especificar'
$ws.Cells.Item(20, 9).Value = 'Primera'
$ws.Cells.Item(20, 10).Value = 20
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 12).Value = 20000
$ws.Cells.Item(20, 13).Value = 20000
$ws.Cells.Item(20, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(20, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(20, 16).Value = 1333
$ws.Cells.Item(20, 17).Value = 15
$ws.Cells.Item(20, 18).Value = 'Hortaliza'
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 21
$ws.Cells.Item(21, 1).Value = 10
$ws.Cells.Item(21, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(21, 3).Value = 'La Araucanía'
$ws.Cells.Item(21, 4).Value = 44179
$ws.Cells.Item(21, 5).Value = 9
$ws.Cells.Item(21, 6).Value = 100114002
$ws.Cells.Item(21, 7).Value = 'Camote'
$ws.Cells.Item(21, 8).Value = 'Sin especificar'
$ws.Cells.Item(21, 9).Value = 'Primera'
$ws.Cells.Item(21, 10).Value = 20
$ws.Cells.Item(21, 11).Value = 20000
$ws.Cells.Item(21, 12).Value = 20000
$ws.Cells.Item(21, 13).Value = 20000
$ws.Cells.Item(21, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(21, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(21, 16).Value = 1333
$ws.Cells.Item(21, 17).Value = 15
$ws.Cells.Item(21, 18).Value = 'Hortaliza'
$ws.Cells.Item(21, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
